$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text ("@") first for cells whose new values would
# otherwise be auto-parsed as numbers by Excel, so they stay inline strings
# matching the source data (e.g. "1.00", "0.610").
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated cell values from the latest crypto price/volume pull.
$ws.Range('D2').Value = '61.803.56'
$ws.Range('E2').Value = '  -2.49%  '
$ws.Range('D3').Value = '2.495.42'
$ws.Range('E3').Value = '  -3.53%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '550.10'
$ws.Range('E5').Value = '  -3.77%  '
$ws.Range('D6').Value = '146.75'
$ws.Range('E6').Value = '  -5.08%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '0.610'
$ws.Range('E8').Value = '  -1.96%  '
$ws.Range('D9').Value = '2.495.77'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('E10').Value = '  -9.24%  '
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').Value = '5.35'
$ws.Range('E12').Value = '  -8.32%  '
$ws.Range('D13').Value = '0.354'
$ws.Range('E13').Value = '  -6.57%  '
$ws.Range('D14').Value = '26.09'
$ws.Range('E14').Value = '  -7.31%  '
$ws.Range('D15').Value = '2.943.65'
$ws.Range('E15').Value = '  -3.70%  '
$ws.Range('D16').Value = '61.711.14'
$ws.Range('E16').Value = '  -2.40%  '
$ws.Range('D17').Value = '0.0000163'
$ws.Range('E17').Value = '  -8.49%  '
$ws.Range('D18').Value = '2.497.11'
$ws.Range('E18').Value = '  -4.56%  '
$ws.Range('D19').Value = '11.10'
$ws.Range('E19').Value = '  -7.00%  '
$ws.Range('D20').Value = '7.01'
$ws.Range('E20').Value = '  -6.30%  '
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  -7.88%  '
$ws.Range('D22').Value = '320.37'
$ws.Range('E22').Value = '  -6.16%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = '63.81'
$ws.Range('E24').Value = '  -5.07%  '
$ws.Range('E25').Value = '  -4.22%  '
$ws.Range('E26').Value = '  -5.77%  '
$ws.Range('D27').Value = '2.626.50'
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').Value = '1.49'
$ws.Range('E29').Value = '  -4.26%  '
$ws.Range('D30').Value = '534.79'
$ws.Range('E30').Value = '  -7.15%  '
$ws.Range('D31').Value = '8.35'
$ws.Range('E31').Value = '  -8.15%  '
$ws.Range('D32').Value = '7.68'
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('E33').Value = '  -6.37%  '
$ws.Range('D34').Value = '1.89'
$ws.Range('E34').Value = '  -7.65%  '
$ws.Range('E35').Value = '  -8.49%  '
$ws.Range('D36').Value = '5.87'
$ws.Range('E36').Value = '  -9.60%  '
$ws.Range('D37').Value = '4.86'
$ws.Range('E37').Value = '  -10.44%  '
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').Value = '0.378'
$ws.Range('E39').Value = '  -5.95%  '
$ws.Range('D40').Value = '18.44'
$ws.Range('E40').Value = '  -6.32%  '
$ws.Range('D41').Value = '143.75'
$ws.Range('E41').Value = '  -6.77%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').Value = '1.69'
$ws.Range('E43').Value = '  -8.71%  '
$ws.Range('D44').Value = '40.38'
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('E45').Value = '  -6.23%  '
$ws.Range('D46').Value = '148.80'
$ws.Range('E46').Value = '  -4.48%  '
$ws.Range('D47').Value = '3.57'
$ws.Range('E47').Value = '  -8.50%  '
$ws.Range('D48').Value = '20.75'
$ws.Range('D49').Value = '0.0534'
$ws.Range('E49').Value = '  -8.88%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.0953'
$ws.Range('E50').Value = '  -4.92%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.590'
$ws.Range('E51').Value = '  -5.54%  '
